$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 62.2
$ws.Range("N2").Value = 85.8724807945396

$ws.Range("K3").Value = 55.8
$ws.Range("N3").Value = 85.8724807945396
